# "new models for renderer"
# Mark rows 6, 7, 9 and 10 (Column D) as "DONE" to match the other
# already-completed rubric rows, and update the sheet view so that the
# window is scrolled back to the top (A1) with G5 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = "DONE"
$ws.Range("D7").Value = "DONE"
$ws.Range("D9").Value = "DONE"
$ws.Range("D10").Value = "DONE"

# Reset scroll position to the top-left of the sheet and select G5,
# matching the new <selection activeCell="G5" sqref="G5"/> view state.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G5").Select()
